$d = $word.ActiveDocument

# The edit removes the blank paragraph plus the two trailing boiler-plate
# paragraphs that were appended right after
# "LOB1012: Estatistica (Requisito)":
#   - (empty paragraph)
#   - "Ver no Jupiter Salvar em pdf Salvar em docx"
#   - "(c) 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github
#      pages. Original theme under Creative Commons Attribution"
# The anchor paragraph itself, and the empty paragraph that follows the
# removed block, are left untouched.

$anchor = "LOB1012: Estat" + [char]0x00ED + "stica (Requisito)"
$expectedAfterAnchor = @(
    "",
    "Ver no Jupiter Salvar em pdf Salvar em docx",
    [char]0x00A9 + " 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution"
)

function Get-ParaText($para) {
    $t = $para.Range.Text
    return $t.TrimEnd([char]0x0007, [char]0x000D, [char]0x000A)
}

# Locate the anchor paragraph.
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ((Get-ParaText $d.Paragraphs.Item($i)) -eq $anchor) {
        $anchorIndex = $i
        break
    }
}

if ($anchorIndex -ge 1) {
    # Delete the paragraph(s) immediately following the anchor, as long as
    # each one's text still matches what we expect to remove. Re-read item
    # ($anchorIndex + 1) each time since prior deletions shift later
    # paragraphs down into that slot.
    foreach ($expected in $expectedAfterAnchor) {
        if (($anchorIndex + 1) -gt $d.Paragraphs.Count) {
            break
        }
        $para = $d.Paragraphs.Item($anchorIndex + 1)
        if ((Get-ParaText $para) -eq $expected) {
            $para.Range.Delete()
        } else {
            break
        }
    }
}
